$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 310.64
$ws.Range("J17").Value = 310.64
$ws.Range("L17").Value = 931.92
$ws.Range("N17").Value = -1267.92

$ws.Range("H28").Value = 2578.3333
$ws.Range("J28").Value = 3050.5
$ws.Range("L28").Value = 3050.5
$ws.Range("N28").Value = -4020.5

$ws.Range("H51").Value = 17748.625
$ws.Range("J51").Value = 8663
$ws.Range("L51").Value = 8663
$ws.Range("N51").Value = -9631

$ws.Range("H92").Value = 64667.562
$ws.Range("I92").Value = 73548.64
$ws.Range("J92").Value = 2500
$ws.Range("K92").Value = 73548.64
$ws.Range("L92").Value = 2500
$ws.Range("M92").Value = -72300.64
$ws.Range("N92").Value = -4996

$ws.Range("H101").Value = 519.8570999999999
$ws.Range("J101").Value = 991.6667
$ws.Range("L101").Value = 2975.0001
$ws.Range("N101").Value = -6219.0001

$ws.Range("H103").Value = 1623.9333
$ws.Range("I103").Value = 1297.5
$ws.Range("J103").Value = 1674.1538
$ws.Range("K103").Value = 3892.5
$ws.Range("L103").Value = 5022.4614
$ws.Range("M103").Value = -3306.5
$ws.Range("N103").Value = -6194.4614

$ws.Range("H137").Value = 2693.6
$ws.Range("I137").Value = 1896.8
$ws.Range("J137").Value = 3490.4
$ws.Range("K137").Value = 5690.4
$ws.Range("L137").Value = 10471.2
$ws.Range("M137").Value = -3140.4
$ws.Range("N137").Value = -15571.2

$ws.Range("H138").Value = 18243.705
$ws.Range("I138").Value = 23858.135
$ws.Range("K138").Value = 71574.405
$ws.Range("M138").Value = -66434.405

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H13").Value = 10858.286
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").ClearContents()

$ws.Range("H97").Value = 804.9048
$ws.Range("I97").Value = 564.9286
$ws.Range("K97").Value = 564.9286
$ws.Range("M97").Value = -68.92859999999996

$ws.Range("H102").Value = 4352.5835
$ws.Range("I102").Value = 4513.1
$ws.Range("J102").Value = 3550
$ws.Range("K102").Value = 4513.1
$ws.Range("L102").Value = 3550
$ws.Range("M102").Value = -2891.1
$ws.Range("N102").Value = -6794

$ws.Range("H105").Value = 100000
$ws.Range("J105").Value = 100000
$ws.Range("L105").Value = 100000
$ws.Range("N105").Value = -106988

$ws.Range("H132").Value = 23332.307
$ws.Range("I132").Value = 27841.275
$ws.Range("J132").Value = 3292.4443
$ws.Range("K132").Value = 83523.82500000001
$ws.Range("L132").Value = 9877.332900000001
$ws.Range("M132").Value = -80993.82500000001
$ws.Range("N132").Value = -14937.3329

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4278.3335
$ws.Range("I86").Value = 3940.25
$ws.Range("J86").Value = 4447.375
$ws.Range("K86").Value = 3940.25
$ws.Range("L86").Value = 4447.375
$ws.Range("M86").Value = -2817.25
$ws.Range("N86").Value = -6693.375

$ws.Range("H89").Value = 4278.3335
$ws.Range("I89").Value = 3940.25
$ws.Range("J89").Value = 4447.375
$ws.Range("K89").Value = 19701.25
$ws.Range("L89").Value = 22236.875
$ws.Range("M89").Value = -14085.25
$ws.Range("N89").Value = -33468.875

$ws.Range("H99").Value = 6256.909
$ws.Range("I99").Value = 5595.1113
$ws.Range("J99").Value = 9235
$ws.Range("K99").Value = 5595.1113
$ws.Range("L99").Value = 9235
$ws.Range("M99").Value = -4097.1113
$ws.Range("N99").Value = -12231

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1653.9333
$ws.Range("I16").Value = 1630.7693
$ws.Range("K16").Value = 1630.7693
$ws.Range("M16").Value = -1343.7693

$ws.Range("H31").Value = 4430.1904
$ws.Range("I31").Value = 2910.6667
$ws.Range("J31").Value = 5569.8335
$ws.Range("K31").Value = 2910.6667
$ws.Range("L31").Value = 5569.8335
$ws.Range("M31").Value = -2615.6667
$ws.Range("N31").Value = -6159.8335

$ws.Range("H34").Value = 4430.1904
$ws.Range("I34").Value = 2910.6667
$ws.Range("J34").Value = 5569.8335
$ws.Range("K34").Value = 2910.6667
$ws.Range("L34").Value = 5569.8335
$ws.Range("M34").Value = -2708.6667
$ws.Range("N34").Value = -5973.8335

$ws.Range("H113").Value = 1653.9333
$ws.Range("I113").Value = 1630.7693
$ws.Range("K113").Value = 1630.7693
$ws.Range("M113").Value = 539.2307000000001

$ws.Range("H134").Value = 67646.625
$ws.Range("I134").Value = 80794.766
$ws.Range("J134").Value = 10671.333
$ws.Range("K134").Value = 242384.298
$ws.Range("L134").Value = 32013.999
$ws.Range("M134").Value = -239849.298
$ws.Range("N134").Value = -37083.999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 621
$ws.Range("I11").Value = 153.6
$ws.Range("K11").Value = 460.8
$ws.Range("M11").Value = -320.8

$ws.Range("H23").Value = 337.25
$ws.Range("J23").Value = 337.25
$ws.Range("L23").Value = 1011.75
$ws.Range("N23").Value = -1481.75

$ws.Range("H69").Value = 1610.3334
$ws.Range("I69").Value = 811.625
$ws.Range("J69").Value = 8000
$ws.Range("K69").Value = 2434.875
$ws.Range("L69").Value = 24000
$ws.Range("M69").Value = -1623.875
$ws.Range("N69").Value = -25622

$ws.Range("H72").Value = 1610.3334
$ws.Range("I72").Value = 811.625
$ws.Range("J72").Value = 8000
$ws.Range("K72").Value = 7304.625
$ws.Range("L72").Value = 72000
$ws.Range("M72").Value = -3248.625
$ws.Range("N72").Value = -80112

$ws.Range("H113").Value = 985
$ws.Range("J113").Value = 898.3333
$ws.Range("L113").Value = 2694.9999
$ws.Range("N113").Value = -7034.9999

$ws.Range("H131").Value = 9858
$ws.Range("I131").Value = 1092.9
$ws.Range("J131").Value = 15336.1875
$ws.Range("K131").Value = 3278.7
$ws.Range("L131").Value = 46008.5625
$ws.Range("M131").Value = 1761.3
$ws.Range("N131").Value = -56088.5625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 5238.316
$ws.Range("I126").Value = 4230.778
$ws.Range("J126").Value = 6145.1
$ws.Range("K126").Value = 12692.334
$ws.Range("L126").Value = 18435.3
$ws.Range("M126").Value = -10222.334
$ws.Range("N126").Value = -23375.3

$ws.Range("H132").Value = 42067
$ws.Range("I132").Value = 54653.473
$ws.Range("J132").Value = 2209.8333
$ws.Range("K132").Value = 163960.419
$ws.Range("L132").Value = 6629.499899999999
$ws.Range("M132").Value = -161430.419
$ws.Range("N132").Value = -11689.4999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2736.8667
$ws.Range("I7").Value = 2218.0715
$ws.Range("K7").Value = 2218.0715
$ws.Range("M7").Value = -2106.0715

$ws.Range("H33").Value = 19999
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()

$ws.Range("H93").Value = 1509.8334
$ws.Range("I93").Value = 1536.9131
$ws.Range("J93").Value = 1420.8572
$ws.Range("K93").Value = 1536.9131
$ws.Range("L93").Value = 1420.8572
$ws.Range("M93").Value = -288.9131
$ws.Range("N93").Value = -3916.8572

$ws.Range("H100").Value = 4179.6
$ws.Range("I100").Value = 4349
$ws.Range("K100").Value = 4349
$ws.Range("M100").Value = -3808

$ws.Range("H126").Value = 2736.8667
$ws.Range("I126").Value = 2218.0715
$ws.Range("K126").Value = 6654.2145
$ws.Range("M126").Value = -4184.2145

$ws.Range("H132").Value = 36721.11
$ws.Range("I132").Value = 50231.96
$ws.Range("J132").Value = 6014.636
$ws.Range("K132").Value = 150695.88
$ws.Range("L132").Value = 18043.908
$ws.Range("M132").Value = -148165.88
$ws.Range("N132").Value = -23103.908

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 18935.666
$ws.Range("J74").Value = 19763
$ws.Range("L74").Value = 19763
$ws.Range("N74").Value = -21635

$ws.Range("H77").Value = 18935.666
$ws.Range("J77").Value = 19763
$ws.Range("L77").Value = 59289
$ws.Range("N77").Value = -68649

$ws.Range("H96").Value = 2631.1667
$ws.Range("I96").Value = 2296
$ws.Range("J96").Value = 2966.3333
$ws.Range("K96").Value = 2296
$ws.Range("L96").Value = 2966.3333
$ws.Range("M96").Value = -923
$ws.Range("N96").Value = -5712.3333

$ws.Range("H100").Value = 1142.1
$ws.Range("I100").Value = 927.75
$ws.Range("J100").Value = 1999.5
$ws.Range("K100").Value = 1855.5
$ws.Range("L100").Value = 3999
$ws.Range("M100").Value = -1314.5
$ws.Range("N100").Value = -5081

$ws.Range("H132").Value = 28035.475
$ws.Range("I132").Value = 32864.688
$ws.Range("J132").Value = 2279.6667
$ws.Range("K132").Value = 98594.06400000001
$ws.Range("L132").Value = 6839.000100000001
$ws.Range("M132").Value = -96064.06400000001
$ws.Range("N132").Value = -11899.0001
